# Generate Report for Handback
# Refresh the generated/handoff/handback timestamps that are written
# whenever the handback status report is (re)generated.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" (row 2 / zh-cn+de-de summary row)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-11-29 04:52:26"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-11-29 04:52:12"
$wsZhCn.Range("K2").Value = "2016-11-29 04:53:04"

# de-de sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-11-29 04:52:26"
$wsDeDe.Range("K2").Value = "2016-11-29 04:53:22"
